$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$h12 = $ws.Range("H12")
try { Write-Host "H12 Font.Color RGB raw:" $h12.Font.Color } catch { Write-Host "ERR" $_.Exception.Message }

# Try using .Interior or OLE color codes
Write-Host "Black RGB should be 0"
